# Generate Report for Handoff
# - Update the "latest status" text and its generated timestamps on all
#   three sheets (Overview, zh-cn, de-de).
# - Narrow the "Status" column(s) that held the long status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Overview sheet --------------------------------------------------
# E2/F2: per-language status; G2: latest handoff xliff generate date
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 00:56:54"

# --- zh-cn sheet -------------------------------------------------------
# C2: Status; H2: Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 00:56:49"

# --- de-de sheet -------------------------------------------------------
# C2: Status; H2: Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 00:56:54"

# --- Column width changes ----------------------------------------------
# Status columns narrowed from the old wide "Handed back..." width down
# to fit "Ready for handoff" (~17.22 chars wide).
$overview.Range("E1").ColumnWidth = 16.333333333333332
$overview.Range("F1").ColumnWidth = 16.333333333333332
$zhcn.Range("C1").ColumnWidth = 16.333333333333332
$dede.Range("C1").ColumnWidth = 16.333333333333332
